$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 72
$ws.Range("I11").Value = 72
$ws.Range("K11").Value = 72
$ws.Range("M11").Value = 68

$ws.Range("H17").Value = 1043.3182
$ws.Range("J17").Value = 1043.3182
$ws.Range("L17").Value = 3129.9546
$ws.Range("N17").Value = -3465.9546

$ws.Range("H29").Value = 1801.5385
$ws.Range("I29").Value = 310.9091
$ws.Range("K29").Value = 932.7273
$ws.Range("M29").Value = -651.7273

$ws.Range("H38").Value = 92.666664
$ws.Range("I38").Value = 92.666664
$ws.Range("K38").Value = 277.999992
$ws.Range("M38").Value = 94.00000799999998

$ws.Range("H51").Value = 26666.334
$ws.Range("J51").Value = 24999.5
$ws.Range("L51").Value = 24999.5
$ws.Range("N51").Value = -25967.5

$ws.Range("H69").Value = 18591.482
$ws.Range("I69").Value = 16664.666
$ws.Range("K69").Value = 49993.99800000001
$ws.Range("M69").Value = -49119.99800000001

$ws.Range("H72").Value = 18591.482
$ws.Range("I72").Value = 16664.666
$ws.Range("K72").Value = 149981.994
$ws.Range("M72").Value = -145613.994

$ws.Range("H76").Value = 168336220
$ws.Range("I76").Value = 336667900
$ws.Range("J76").Value = 4555
$ws.Range("K76").Value = 336667900
$ws.Range("L76").Value = 4555
$ws.Range("M76").Value = -336667585
$ws.Range("N76").Value = -5185

$ws.Range("H79").Value = 168336220
$ws.Range("I79").Value = 336667900
$ws.Range("J79").Value = 4555
$ws.Range("K79").Value = 336667900
$ws.Range("L79").Value = 4555
$ws.Range("M79").Value = -336666808
$ws.Range("N79").Value = -6739

$ws.Range("H98").Value = 591.8
$ws.Range("I98").Value = 591.8
$ws.Range("K98").Value = 591.8
$ws.Range("M98").Value = 906.2

$ws.Range("H116").Value = 13523.667
$ws.Range("I116").Value = 13127
$ws.Range("J116").Value = 13841
$ws.Range("K116").Value = 13127
$ws.Range("L116").Value = 13841
$ws.Range("M116").Value = -9685
$ws.Range("N116").Value = -20725

$ws.Range("H122").Value = 591.8
$ws.Range("I122").Value = 591.8
$ws.Range("K122").Value = 1775.4
$ws.Range("M122").Value = 674.6000000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2472.889
$ws.Range("I45").Value = 1717
$ws.Range("K45").Value = 1717
$ws.Range("M45").Value = -1340

$ws.Range("H61").Value = 12818.333
$ws.Range("I61").Value = 10387.792
$ws.Range("J61").Value = 19299.777
$ws.Range("K61").Value = 10387.792
$ws.Range("L61").Value = 19299.777
$ws.Range("M61").Value = -10175.792
$ws.Range("N61").Value = -19723.777

$ws.Range("H102").Value = 1109.8334
$ws.Range("I102").Value = 1157.5883
$ws.Range("J102").Value = 298
$ws.Range("K102").Value = 1157.5883
$ws.Range("L102").Value = 298
$ws.Range("M102").Value = 464.4117000000001
$ws.Range("N102").Value = -3542

$ws.Range("H136").Value = 12818.333
$ws.Range("I136").Value = 10387.792
$ws.Range("J136").Value = 19299.777
$ws.Range("K136").Value = 31163.376
$ws.Range("L136").Value = 57899.33099999999
$ws.Range("M136").Value = -28613.376
$ws.Range("N136").Value = -62999.33099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H86").Value = 11660.857
$ws.Range("J86").Value = 23460.166
$ws.Range("L86").Value = 23460.166
$ws.Range("N86").Value = -25706.166

$ws.Range("H89").Value = 11660.857
$ws.Range("J89").Value = 23460.166
$ws.Range("L89").Value = 117300.83
$ws.Range("N89").Value = -128532.83

$ws.Range("H105").Value = 3716.9167
$ws.Range("I105").Value = 3632.4285
$ws.Range("K105").Value = 3632.4285
$ws.Range("M105").Value = -1885.4285

$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178

$ws.Range("H134").Value = 6140.6606
$ws.Range("I134").Value = 4941.75
$ws.Range("K134").Value = 14825.25
$ws.Range("M134").Value = -12290.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2695.825
$ws.Range("I31").Value = 1719.9
$ws.Range("J31").Value = 3021.1333
$ws.Range("K31").Value = 1719.9
$ws.Range("L31").Value = 3021.1333
$ws.Range("M31").Value = -1424.9
$ws.Range("N31").Value = -3611.1333

$ws.Range("H34").Value = 2695.825
$ws.Range("I34").Value = 1719.9
$ws.Range("J34").Value = 3021.1333
$ws.Range("K34").Value = 1719.9
$ws.Range("L34").Value = 3021.1333
$ws.Range("M34").Value = -1517.9
$ws.Range("N34").Value = -3425.1333

$ws.Range("H58").Value = 3176.7693
$ws.Range("I58").Value = 2399.16
$ws.Range("K58").Value = 2399.16
$ws.Range("M58").Value = -2196.16

$ws.Range("H107").Value = 407.33334
$ws.Range("I107").Value = 345.6154
$ws.Range("K107").Value = 345.6154
$ws.Range("M107").Value = 1574.3846

$ws.Range("H120").Value = 97326
$ws.Range("J120").Value = 97326
$ws.Range("L120").Value = 97326
$ws.Range("N120").Value = -104584

$ws.Range("H132").Value = 2167.7036
$ws.Range("I132").Value = 1842.2273
$ws.Range("K132").Value = 5526.6819
$ws.Range("M132").Value = -2996.6819

$ws.Range("H134").Value = 4156.2373
$ws.Range("I134").Value = 3708.2246
$ws.Range("J134").Value = 9140.375
$ws.Range("K134").Value = 11124.6738
$ws.Range("L134").Value = 27421.125
$ws.Range("M134").Value = -8589.6738
$ws.Range("N134").Value = -32491.125

$ws.Range("H136").Value = 3176.7693
$ws.Range("I136").Value = 2399.16
$ws.Range("K136").Value = 7197.48
$ws.Range("M136").Value = -4647.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 2061.8
$ws.Range("I58").Value = 1766.3334
$ws.Range("K58").Value = 5299.0002
$ws.Range("M58").Value = -5171.0002

$ws.Range("H62").Value = 9312.375
$ws.Range("J62").Value = 9999.857
$ws.Range("L62").Value = 29999.571
$ws.Range("N62").Value = -31371.571

$ws.Range("H65").Value = 9312.375
$ws.Range("J65").Value = 9999.857
$ws.Range("L65").Value = 89998.713
$ws.Range("N65").Value = -96862.713

$ws.Range("H107").Value = 579.8261
$ws.Range("I107").Value = 256.57144
$ws.Range("J107").Value = 721.25
$ws.Range("K107").Value = 769.71432
$ws.Range("L107").Value = 2163.75
$ws.Range("M107").Value = 1150.28568
$ws.Range("N107").Value = -6003.75

$ws.Range("H116").Value = 50249.5
$ws.Range("I116").Value = 50249.5
$ws.Range("K116").Value = 150748.5
$ws.Range("M116").Value = -147306.5

$ws.Range("H121").Value = 54792.086
$ws.Range("I121").Value = 886
$ws.Range("J121").Value = 69766
$ws.Range("K121").Value = 2658
$ws.Range("L121").Value = 209298
$ws.Range("M121").Value = -1348
$ws.Range("N121").Value = -211918

$ws.Range("H129").Value = 9260942
$ws.Range("I129").Value = 446.7
$ws.Range("K129").Value = 1340.1
$ws.Range("M129").Value = 3659.9

$ws.Range("H136").Value = 18477
$ws.Range("I136").Value = 6596.7144
$ws.Range("J136").Value = 35109.4
$ws.Range("K136").Value = 19790.1432
$ws.Range("L136").Value = 105328.2
$ws.Range("M136").Value = -14690.1432
$ws.Range("N136").Value = -115528.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 270684
$ws.Range("J108").Value = 270684
$ws.Range("L108").Value = 270684
$ws.Range("N108").Value = -278364

$ws.Range("H132").Value = 1652.5
$ws.Range("I132").Value = 1129.3334
$ws.Range("J132").Value = 2437.25
$ws.Range("K132").Value = 3388.0002
$ws.Range("L132").Value = 7311.75
$ws.Range("M132").Value = -858.0001999999999
$ws.Range("N132").Value = -12371.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1962.5454
$ws.Range("I22").Value = 1758
$ws.Range("J22").Value = 2022.7059
$ws.Range("K22").Value = 1758
$ws.Range("L22").Value = 2022.7059
$ws.Range("M22").Value = -1463
$ws.Range("N22").Value = -2612.7059

$ws.Range("H27").Value = 1962.5454
$ws.Range("I27").Value = 1758
$ws.Range("J27").Value = 2022.7059
$ws.Range("K27").Value = 1758
$ws.Range("L27").Value = 2022.7059
$ws.Range("M27").Value = -1651
$ws.Range("N27").Value = -2236.7059

$ws.Range("H46").Value = 1795.65
$ws.Range("I46").Value = 1179.5714
$ws.Range("K46").Value = 1179.5714
$ws.Range("M46").Value = -991.5714

$ws.Range("H55").Value = 1325
$ws.Range("I55").Value = 300
$ws.Range("J55").Value = 1581.25
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 1581.25
$ws.Range("M55").Value = -127
$ws.Range("N55").Value = -1927.25

$ws.Range("H68").Value = 2633.4614
$ws.Range("I68").Value = 2529.375
$ws.Range("J68").Value = 2800
$ws.Range("K68").Value = 2529.375
$ws.Range("L68").Value = 2800
$ws.Range("M68").Value = -1780.375
$ws.Range("N68").Value = -4298

$ws.Range("H71").Value = 2633.4614
$ws.Range("I71").Value = 2529.375
$ws.Range("J71").Value = 2800
$ws.Range("K71").Value = 12646.875
$ws.Range("L71").Value = 14000
$ws.Range("M71").Value = -8902.875
$ws.Range("N71").Value = -21488

$ws.Range("H132").Value = 2709.5806
$ws.Range("I132").Value = 1702.091
$ws.Range("J132").Value = 3263.7
$ws.Range("K132").Value = 5106.272999999999
$ws.Range("L132").Value = 9791.099999999999
$ws.Range("M132").Value = -2576.272999999999
$ws.Range("N132").Value = -14851.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 822.4815
$ws.Range("I113").Value = 525.1539
$ws.Range("J113").Value = 1098.5714
$ws.Range("K113").Value = 1575.4617
$ws.Range("L113").Value = 3295.7142
$ws.Range("M113").Value = 594.5382999999999
$ws.Range("N113").Value = -7635.7142
